$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Nudge "Straight Connector 66" (id=67) slightly to the right:
#    a:off x 7834149 -> 7843674 (y, cx, cy unchanged)
$conn66 = $s.Shapes.Item("Straight Connector 66")
$conn66.Left = 617.6122

# 2) Remove "Straight Connector 90" (id=91)
$conn90 = $s.Shapes.Item("Straight Connector 90")
$conn90.Delete()

# 3) Add a new straight connector/arrow ("Straight Connector 1") at the end
#    of the shape tree, duplicating the style of the nearby connector so it
#    keeps the matching line/arrowhead formatting and shape style refs.
$srcConn = $s.Shapes.Item("Straight Connector 66")
$dup = $srcConn.Duplicate()
$newConn = $dup.Item(1)
$newConn.Name = "Straight Connector 1"
$newConn.Left = 780.7396850393701
$newConn.Top = 79.5757485
$newConn.Width = 28.3464575
$newConn.Height = 0
